$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44308, 55, 255, 134.9113553036035),
    @(44309, 32, 260, 137.556675995831),
    @(44310, 39, 248, 131.2079063344849),
    @(44311, 39, 252, 133.324162888267),
    @(44312, 61, 254, 134.3822911651579)
)

$startRow = 234
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Range("A233:D233").Copy()
    $ws.Range(("A" + $r + ":D" + $r)).PasteSpecial(-4122)
    $row = $data[$i]
    $ws.Range("A" + $r).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]
}
